$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 36

# Trade # (number)
$ws.Cells.Item($row, 1).Value = 46

# Date - stored as plain text in this sheet (e.g. "2026-02-16"), not a date
# serial. Force a text number format before assigning so the COM layer does
# not auto-convert the string into a date value, then clear the format again
# so no stray style is left behind (matching the rest of the sheet, which
# has no per-cell styles).
$dateCell = $ws.Cells.Item($row, 2)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-02-16"
$dateCell.ClearFormats()

# Time (plain text, not auto-converted by the engine)
$ws.Cells.Item($row, 3).Value = "21:30:05"

# Strategy / Side
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "DOWN"

# Entry price
$ws.Cells.Item($row, 6).Value = 68641.42999999999

# Exit price - left blank (open trade), matches the empty "Exit Price" cells
# used throughout this sheet for still-open trades.

# Status
$ws.Cells.Item($row, 8).Value = "OPEN"

# P&L % / P&L $
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0

# Confidence
$ws.Cells.Item($row, 11).Value = 0.75

# Entry reason
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.163% move"

# Exit reason - left blank (open trade)

# Duration (min)
$ws.Cells.Item($row, 14).Value = 0
